$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.035.35"
$ws.Range("E2").Value = "  +4.39%  "
$ws.Range("D3").Value = "3.676.23"
$ws.Range("E3").Value = "  +10.48%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.93"
$ws.Range("E5").Value = "  +5.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "645.92"
$ws.Range("E6").Value = "  +5.25%  "
$ws.Range("E7").Value = "  +5.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.02"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("D11").Value = "3.675.62"
$ws.Range("E11").Value = "  +10.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.06"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.40"
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("D15").Value = "4.373.93"
$ws.Range("E15").Value = "  +10.72%  "
$ws.Range("D16").Value = "95.868.08"
$ws.Range("E16").Value = "  +4.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000258"
$ws.Range("E17").Value = "  +6.08%  "
$ws.Range("D18").Value = "3.677.81"
$ws.Range("E18").Value = "  +10.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.45"
$ws.Range("E19").Value = "  +24.16%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.74"
$ws.Range("E21").Value = "  +8.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "521.68"
$ws.Range("E22").Value = "  +6.12%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.483"
$ws.Range("E24").Value = "  +10.50%  "
$ws.Range("E25").Value = "  +9.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.84"
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "97.54"
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.70"
$ws.Range("E28").Value = "  +6.64%  "
$ws.Range("E29").Value = "  +21.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.75"
$ws.Range("E30").Value = "  +5.73%  "
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.18"
$ws.Range("E33").Value = "  +17.36%  "
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.181"
$ws.Range("E34").Value = "  +4.37%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "563.35"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.98"
$ws.Range("E38").Value = "  +7.53%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +10.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.970"
$ws.Range("E40").Value = "  +11.65%  "
$ws.Range("E41").Value = "  +3.03%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.83"
$ws.Range("E44").Value = "  +7.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0431"
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.71"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.09"
$ws.Range("E47").Value = "  +51.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  +5.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.91"
$ws.Range("E49").Value = "  +5.50%  "
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("E51").Value = "  -2.82%  "
